$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 145.5
$ws.Range("I9").Value = 136.5
$ws.Range("J9").Value = 199.5
$ws.Range("K9").Value = 136.5
$ws.Range("L9").Value = 199.5
$ws.Range("M9").Value = 32.5
$ws.Range("N9").Value = -537.5
$ws.Range("H13").Value = 15000
$ws.Range("I13").Value = 15000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 15000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -14831
$ws.Range("N13").ClearContents()
$ws.Range("H15").Value = 1332.4286
$ws.Range("I15").Value = 1332.4286
$ws.Range("K15").Value = 3997.2858
$ws.Range("M15").Value = -3828.2858
$ws.Range("H16").Value = 63333.332
$ws.Range("J16").Value = 63333.332
$ws.Range("L16").Value = 63333.332
$ws.Range("N16").Value = -63793.332
$ws.Range("H18").Value = 500
$ws.Range("I18").Value = 500
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 500
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -216
$ws.Range("N18").Value = -1068
$ws.Range("H33").Value = 746.4211
$ws.Range("I33").Value = 460
$ws.Range("K33").Value = 460
$ws.Range("M33").Value = -231
$ws.Range("H99").Value = 492.8
$ws.Range("J99").Value = 750
$ws.Range("L99").Value = 2250
$ws.Range("N99").Value = -5246
$ws.Range("H100").Value = 2839.4707
$ws.Range("I100").Value = 2436.1
$ws.Range("J100").Value = 3415.7144
$ws.Range("K100").Value = 2436.1
$ws.Range("L100").Value = 3415.7144
$ws.Range("M100").Value = -1895.1
$ws.Range("N100").Value = -4497.7144
$ws.Range("H101").Value = 451.5
$ws.Range("I101").Value = 451.5
$ws.Range("K101").Value = 1354.5
$ws.Range("M101").Value = 267.5
$ws.Range("H106").Value = 26773.066
$ws.Range("I106").Value = 26773.066
$ws.Range("K106").Value = 26773.066
$ws.Range("M106").Value = -26142.066
$ws.Range("H111").Value = 3110
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()
$ws.Range("H112").Value = 2324.5
$ws.Range("J112").Value = 2355.9167
$ws.Range("L112").Value = 7067.750100000001
$ws.Range("N112").Value = -9283.750100000001
$ws.Range("H115").Value = 7262.1
$ws.Range("I115").Value = 3919.3333
$ws.Range("J115").Value = 8694.714
$ws.Range("K115").Value = 11757.9999
$ws.Range("L115").Value = 26084.142
$ws.Range("M115").Value = -10190.9999
$ws.Range("N115").Value = -29218.142
$ws.Range("H118").Value = 4507.0835
$ws.Range("I118").Value = 4528.9
$ws.Range("K118").Value = 13586.7
$ws.Range("M118").Value = -11929.7
$ws.Range("H121").Value = 1997.909
$ws.Range("J121").Value = 1997.909
$ws.Range("L121").Value = 5993.727000000001
$ws.Range("N121").Value = -9487.727000000001
$ws.Range("H124").Value = 99995
$ws.Range("J124").Value = 99995
$ws.Range("L124").Value = 99995
$ws.Range("N124").Value = -109815
$ws.Range("H129").Value = 2288.9285
$ws.Range("J129").Value = 3162.6667
$ws.Range("L129").Value = 9488.000100000001
$ws.Range("N129").Value = -19488.0001
$ws.Range("H130").Value = 87985
$ws.Range("J130").Value = 87985
$ws.Range("L130").Value = 87985
$ws.Range("N130").Value = -98025
$ws.Range("H131").Value = 19729.666
$ws.Range("I131").Value = 12944.5
$ws.Range("K131").Value = 38833.5
$ws.Range("M131").Value = -33793.5
$ws.Range("H132").Value = 5204.107
$ws.Range("I132").Value = 2828.84
$ws.Range("K132").Value = 8486.52
$ws.Range("M132").Value = -5956.52
$ws.Range("H133").Value = 99995
$ws.Range("J133").Value = 99995
$ws.Range("L133").Value = 99995
$ws.Range("N133").Value = -110115
$ws.Range("H135").Value = 4164.4443
$ws.Range("I135").Value = 4164.4443
$ws.Range("K135").Value = 37479.9987
$ws.Range("M135").Value = -34944.9987
$ws.Range("H138").Value = 2296.7234
$ws.Range("I138").Value = 1998.6666
$ws.Range("J138").Value = 2317.0454
$ws.Range("K138").Value = 5995.9998
$ws.Range("L138").Value = 6951.1362
$ws.Range("M138").Value = -855.9997999999996
$ws.Range("N138").Value = -17231.1362
$ws.Range("H139").Value = 99995
$ws.Range("J139").Value = 99995
$ws.Range("L139").Value = 99995
$ws.Range("N139").Value = -110275

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3254
$ws.Range("I2").Value = 1698.9231
$ws.Range("K2").Value = 1698.9231
$ws.Range("M2").Value = -1585.9231
$ws.Range("H4").Value = 594.2
$ws.Range("I4").Value = 568
$ws.Range("K4").Value = 568
$ws.Range("M4").Value = -452
$ws.Range("H5").Value = 2956.7778
$ws.Range("I5").Value = 903.1667
$ws.Range("J5").Value = 7064
$ws.Range("K5").Value = 903.1667
$ws.Range("L5").Value = 7064
$ws.Range("M5").Value = -791.1667
$ws.Range("N5").Value = -7288
$ws.Range("H12").Value = 1600
$ws.Range("I12").Value = 2000
$ws.Range("J12").Value = 1200
$ws.Range("K12").Value = 2000
$ws.Range("L12").Value = 1200
$ws.Range("M12").Value = -1827
$ws.Range("N12").Value = -1546
$ws.Range("H32").Value = 1662.4938
$ws.Range("I32").Value = 1180.4474
$ws.Range("K32").Value = 1180.4474
$ws.Range("M32").Value = -893.4474
$ws.Range("H49").Value = 69865
$ws.Range("J49").Value = 69865
$ws.Range("L49").Value = 69865
$ws.Range("N49").Value = -70385
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H58").Value = 99533.5
$ws.Range("J58").Value = 99533.5
$ws.Range("L58").Value = 99533.5
$ws.Range("N58").Value = -100393.5
$ws.Range("H61").Value = 5824.7144
$ws.Range("I61").Value = 5489.6665
$ws.Range("K61").Value = 5489.6665
$ws.Range("M61").Value = -5277.6665
$ws.Range("H63").Value = 1192.3334
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1192.3334
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H92").Value = 49850
$ws.Range("J92").Value = 49850
$ws.Range("L92").Value = 49850
$ws.Range("N92").Value = -54842
$ws.Range("H97").Value = 771.1667
$ws.Range("J97").Value = 1449.75
$ws.Range("L97").Value = 1449.75
$ws.Range("N97").Value = -2441.75
$ws.Range("H101").Value = 278200.6
$ws.Range("J101").Value = 278200.6
$ws.Range("L101").Value = 278200.6
$ws.Range("N101").Value = -284690.6
$ws.Range("H102").Value = 6140.778
$ws.Range("I102").Value = 6140.778
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 6140.778
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -4518.778
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 3254
$ws.Range("I116").Value = 1698.9231
$ws.Range("K116").Value = 1698.9231
$ws.Range("M116").Value = 595.0769
$ws.Range("H122").Value = 2547.5405
$ws.Range("I122").Value = 3068.8635
$ws.Range("K122").Value = 9206.5905
$ws.Range("M122").Value = -6756.5905
$ws.Range("H123").Value = 35424.5
$ws.Range("J123").Value = 35424.5
$ws.Range("L123").Value = 35424.5
$ws.Range("N123").Value = -45224.5
$ws.Range("H132").Value = 4043.8215
$ws.Range("I132").Value = 4596.303
$ws.Range("K132").Value = 13788.909
$ws.Range("M132").Value = -11258.909
$ws.Range("H136").Value = 5824.7144
$ws.Range("I136").Value = 5489.6665
$ws.Range("K136").Value = 16468.9995
$ws.Range("M136").Value = -13918.9995
$ws.Range("H141").Value = 100000
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 100000
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 100000
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -110360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3254
$ws.Range("I3").Value = 1698.9231
$ws.Range("K3").Value = 1698.9231
$ws.Range("M3").Value = -1584.9231
$ws.Range("H4").Value = 2956.7778
$ws.Range("I4").Value = 903.1667
$ws.Range("J4").Value = 7064
$ws.Range("K4").Value = 903.1667
$ws.Range("L4").Value = 7064
$ws.Range("M4").Value = -788.1667
$ws.Range("N4").Value = -7294
$ws.Range("H22").Value = 649.3333
$ws.Range("I22").Value = 649.3333
$ws.Range("K22").Value = 649.3333
$ws.Range("M22").Value = -476.3333
$ws.Range("H94").Value = 1267.8125
$ws.Range("I94").Value = 1320.3572
$ws.Range("K94").Value = 1320.3572
$ws.Range("M94").Value = -869.3571999999999
$ws.Range("H99").Value = 2848.8462
$ws.Range("I99").Value = 3381.5
$ws.Range("J99").Value = 1073.3334
$ws.Range("K99").Value = 3381.5
$ws.Range("L99").Value = 1073.3334
$ws.Range("M99").Value = -1883.5
$ws.Range("N99").Value = -4069.3334
$ws.Range("H105").Value = 1868.6
$ws.Range("I105").Value = 1197.4
$ws.Range("J105").Value = 2539.8
$ws.Range("K105").Value = 1197.4
$ws.Range("L105").Value = 2539.8
$ws.Range("M105").Value = 549.5999999999999
$ws.Range("N105").Value = -6033.8
$ws.Range("H106").Value = 12667.333
$ws.Range("J106").Value = 12667.333
$ws.Range("L106").Value = 12667.333
$ws.Range("N106").Value = -15191.333
$ws.Range("H134").Value = 11105.833
$ws.Range("I134").Value = 14579.625
$ws.Range("J134").Value = 8326.799999999999
$ws.Range("K134").Value = 43738.875
$ws.Range("L134").Value = 24980.4
$ws.Range("M134").Value = -41203.875
$ws.Range("N134").Value = -30050.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 7584
$ws.Range("I3").Value = 4682.6665
$ws.Range("K3").Value = 4682.6665
$ws.Range("M3").Value = -4569.6665
$ws.Range("H7").Value = 151.25
$ws.Range("I7").Value = 141.75
$ws.Range("J7").Value = 170.25
$ws.Range("K7").Value = 141.75
$ws.Range("L7").Value = 170.25
$ws.Range("M7").Value = -28.75
$ws.Range("N7").Value = -396.25
$ws.Range("H15").Value = 14999
$ws.Range("J15").Value = 14999
$ws.Range("L15").Value = 14999
$ws.Range("N15").Value = -15339
$ws.Range("H16").Value = 916.8570999999999
$ws.Range("I16").Value = 774.7692
$ws.Range("J16").Value = 1147.75
$ws.Range("K16").Value = 774.7692
$ws.Range("L16").Value = 1147.75
$ws.Range("M16").Value = -487.7692
$ws.Range("N16").Value = -1721.75
$ws.Range("H22").Value = 4500
$ws.Range("I22").Value = 5333.3335
$ws.Range("K22").Value = 5333.3335
$ws.Range("M22").Value = -4983.3335
$ws.Range("H25").Value = 9997.5
$ws.Range("J25").Value = 9997
$ws.Range("L25").Value = 9997
$ws.Range("N25").Value = -10345
$ws.Range("H43").Value = 26412
$ws.Range("J43").Value = 26412
$ws.Range("L43").Value = 26412
$ws.Range("N43").Value = -26780
$ws.Range("H88").Value = 27499.334
$ws.Range("J88").Value = 24749
$ws.Range("L88").Value = 24749
$ws.Range("N88").Value = -25561
$ws.Range("H91").Value = 27499.334
$ws.Range("J91").Value = 24749
$ws.Range("L91").Value = 24749
$ws.Range("N91").Value = -27557
$ws.Range("H95").Value = 22512.334
$ws.Range("J95").Value = 22512.334
$ws.Range("L95").Value = 22512.334
$ws.Range("N95").Value = -28004.334
$ws.Range("H96").Value = 16158
$ws.Range("J96").Value = 16158
$ws.Range("L96").Value = 16158
$ws.Range("N96").Value = -21650
$ws.Range("H99").Value = 2346
$ws.Range("I99").Value = 2012
$ws.Range("J99").Value = 3014
$ws.Range("K99").Value = 2012
$ws.Range("L99").Value = 3014
$ws.Range("M99").Value = -514
$ws.Range("N99").Value = -6010
$ws.Range("H101").Value = 26412
$ws.Range("J101").Value = 26412
$ws.Range("L101").Value = 26412
$ws.Range("N101").Value = -32902
$ws.Range("H105").Value = 1494.0333
$ws.Range("I105").Value = 1595.3572
$ws.Range("J105").Value = 1405.375
$ws.Range("K105").Value = 1595.3572
$ws.Range("L105").Value = 1405.375
$ws.Range("M105").Value = 151.6428000000001
$ws.Range("N105").Value = -4899.375
$ws.Range("H107").Value = 1015.1539
$ws.Range("J107").Value = 1147.75
$ws.Range("L107").Value = 1147.75
$ws.Range("N107").Value = -4987.75
$ws.Range("H113").Value = 916.8570999999999
$ws.Range("I113").Value = 774.7692
$ws.Range("J113").Value = 1147.75
$ws.Range("K113").Value = 774.7692
$ws.Range("L113").Value = 1147.75
$ws.Range("M113").Value = 1395.2308
$ws.Range("N113").Value = -5487.75
$ws.Range("H126").Value = 2346
$ws.Range("I126").Value = 2012
$ws.Range("J126").Value = 3014
$ws.Range("K126").Value = 6036
$ws.Range("L126").Value = 9042
$ws.Range("M126").Value = -3566
$ws.Range("N126").Value = -13982

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 555613.9399999999
$ws.Range("I2").Value = 800030.5600000001
$ws.Range("J2").Value = 121.545456
$ws.Range("K2").Value = 4800183.36
$ws.Range("L2").Value = 729.272736
$ws.Range("M2").Value = -4800070.36
$ws.Range("N2").Value = -955.272736
$ws.Range("H12").Value = 259.7143
$ws.Range("I12").Value = 401
$ws.Range("J12").Value = 153.75
$ws.Range("K12").Value = 1203
$ws.Range("L12").Value = 461.25
$ws.Range("M12").Value = -1030
$ws.Range("N12").Value = -807.25
$ws.Range("H17").Value = 1215.4286
$ws.Range("I17").Value = 337
$ws.Range("J17").Value = 1874.25
$ws.Range("K17").Value = 1011
$ws.Range("L17").Value = 5622.75
$ws.Range("M17").Value = -842
$ws.Range("N17").Value = -5960.75
$ws.Range("H46").Value = 2500518.5
$ws.Range("J46").Value = 5000599.5
$ws.Range("L46").Value = 15001798.5
$ws.Range("N46").Value = -15001980.5
$ws.Range("H51").Value = 1026.7142
$ws.Range("I51").Value = 422.75
$ws.Range("K51").Value = 1268.25
$ws.Range("M51").Value = -808.25
$ws.Range("H68").Value = 1683.7
$ws.Range("J68").Value = 2119.889
$ws.Range("L68").Value = 6359.667
$ws.Range("N68").Value = -7981.667
$ws.Range("H71").Value = 1683.7
$ws.Range("J71").Value = 2119.889
$ws.Range("L71").Value = 19079.001
$ws.Range("N71").Value = -27191.001
$ws.Range("H99").Value = 5766.3125
$ws.Range("I99").Value = 2932.818
$ws.Range("K99").Value = 8798.454000000002
$ws.Range("M99").Value = -6552.454000000002
$ws.Range("H107").Value = 981.08
$ws.Range("J107").Value = 1236.6364
$ws.Range("L107").Value = 3709.9092
$ws.Range("N107").Value = -7549.9092
$ws.Range("H121").Value = 1347.3871
$ws.Range("I121").Value = 1000
$ws.Range("K121").Value = 3000
$ws.Range("M121").Value = -1690
$ws.Range("H132").Value = 1910.4445
$ws.Range("I132").Value = 1469.8572
$ws.Range("J132").Value = 3452.5
$ws.Range("K132").Value = 13228.7148
$ws.Range("L132").Value = 31072.5
$ws.Range("M132").Value = -10698.7148
$ws.Range("N132").Value = -36132.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 14370.444
$ws.Range("J24").Value = 14370.444
$ws.Range("L24").Value = 14370.444
$ws.Range("N24").Value = -14716.444
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("N50").ClearContents()
$ws.Range("H80").Value = 1424.5
$ws.Range("I80").Value = 1424.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1424.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -426.5
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 1424.5
$ws.Range("I83").Value = 1424.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7122.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2130.5
$ws.Range("N83").ClearContents()
$ws.Range("H98").Value = 19046.5
$ws.Range("J98").Value = 19046.5
$ws.Range("L98").Value = 19046.5
$ws.Range("N98").Value = -25036.5
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").ClearContents()
$ws.Range("H132").Value = 8959.630999999999
$ws.Range("I132").Value = 8553.714
$ws.Range("J132").Value = 10251.182
$ws.Range("K132").Value = 25661.142
$ws.Range("L132").Value = 30753.546
$ws.Range("M132").Value = -23131.142
$ws.Range("N132").Value = -35813.546

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 84
$ws.Range("I9").Value = 84
$ws.Range("K9").Value = 84
$ws.Range("M9").Value = 140
$ws.Range("H11").Value = 4999.5
$ws.Range("I11").Value = 5000
$ws.Range("J11").Value = 4999
$ws.Range("K11").Value = 5000
$ws.Range("L11").Value = 4999
$ws.Range("M11").Value = -4860
$ws.Range("N11").Value = -5279
$ws.Range("H16").Value = 730.4483
$ws.Range("I16").Value = 563.7143
$ws.Range("K16").Value = 563.7143
$ws.Range("M16").Value = -393.7143
$ws.Range("H22").Value = 2150.8823
$ws.Range("J22").Value = 2150.8823
$ws.Range("L22").Value = 2150.8823
$ws.Range("N22").Value = -2740.8823
$ws.Range("H27").Value = 2150.8823
$ws.Range("J27").Value = 2150.8823
$ws.Range("L27").Value = 2150.8823
$ws.Range("N27").Value = -2364.8823
$ws.Range("H36").Value = 38666
$ws.Range("J36").Value = 38666
$ws.Range("L36").Value = 38666
$ws.Range("N36").Value = -39790
$ws.Range("H40").Value = 3673.2903
$ws.Range("I40").Value = 3832.1482
$ws.Range("J40").Value = 2601
$ws.Range("K40").Value = 3832.1482
$ws.Range("L40").Value = 2601
$ws.Range("M40").Value = -3696.1482
$ws.Range("N40").Value = -2873
$ws.Range("H48").Value = 28000
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 28000
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 28000
$ws.Range("M48").ClearContents()
$ws.Range("N48").Value = -29322
$ws.Range("H55").Value = 1025.3182
$ws.Range("I55").Value = 547.2
$ws.Range("J55").Value = 1165.9412
$ws.Range("K55").Value = 547.2
$ws.Range("L55").Value = 1165.9412
$ws.Range("M55").Value = -374.2
$ws.Range("N55").Value = -1511.9412
$ws.Range("H68").Value = 5198.615
$ws.Range("I68").Value = 3311
$ws.Range("K68").Value = 3311
$ws.Range("M68").Value = -2562
$ws.Range("H71").Value = 5198.615
$ws.Range("I71").Value = 3311
$ws.Range("K71").Value = 16555
$ws.Range("M71").Value = -12811
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H100").Value = 1786.7084
$ws.Range("I100").Value = 1447.238
$ws.Range("K100").Value = 1447.238
$ws.Range("M100").Value = -906.2380000000001
$ws.Range("H106").Value = 8992.5
$ws.Range("J106").Value = 8992.5
$ws.Range("L106").Value = 8992.5
$ws.Range("N106").Value = -11516.5
$ws.Range("H122").Value = 4634.9165
$ws.Range("I122").Value = 4656.4546
$ws.Range("J122").Value = 4616.6924
$ws.Range("K122").Value = 13969.3638
$ws.Range("L122").Value = 13850.0772
$ws.Range("M122").Value = -11519.3638
$ws.Range("N122").Value = -18750.0772
$ws.Range("H124").Value = 49999
$ws.Range("J124").Value = 49999
$ws.Range("L124").Value = 49999
$ws.Range("N124").Value = -59819
$ws.Range("H132").Value = 8650.366
$ws.Range("I132").Value = 8309.620999999999
$ws.Range("J132").Value = 9473.833000000001
$ws.Range("K132").Value = 24928.863
$ws.Range("L132").Value = 28421.499
$ws.Range("M132").Value = -22398.863
$ws.Range("N132").Value = -33481.499
$ws.Range("H136").Value = 2706.7058
$ws.Range("I136").Value = 2052.6428
$ws.Range("K136").Value = 6157.928400000001
$ws.Range("M136").Value = -3607.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 15000
$ws.Range("I31").Value = 15000
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 15000
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -14652
$ws.Range("N31").ClearContents()
$ws.Range("H34").Value = 69999.336
$ws.Range("I34").Value = 49999
$ws.Range("J34").Value = 79999.5
$ws.Range("K34").Value = 49999
$ws.Range("L34").Value = 79999.5
$ws.Range("M34").Value = -49796
$ws.Range("N34").Value = -80405.5
$ws.Range("H37").Value = 49999
$ws.Range("I37").Value = 49999
$ws.Range("K37").Value = 49999
$ws.Range("M37").Value = -49796
$ws.Range("H42").Value = 49249.5
$ws.Range("J42").Value = 49000
$ws.Range("L42").Value = 49000
$ws.Range("N42").Value = -49756
$ws.Range("H81").Value = 1640.8572
$ws.Range("I81").Value = 1640.8572
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3281.7144
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2220.7144
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 1640.8572
$ws.Range("I84").Value = 1640.8572
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 16408.572
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -11104.572
$ws.Range("N84").ClearContents()
$ws.Range("H103").Value = 25867
$ws.Range("J103").Value = 25867
$ws.Range("L103").Value = 25867
$ws.Range("N103").Value = -28211
$ws.Range("H105").Value = 66648.60000000001
$ws.Range("J105").Value = 66648.60000000001
$ws.Range("L105").Value = 66648.60000000001
$ws.Range("N105").Value = -73636.60000000001
$ws.Range("H132").Value = 2929.0789
$ws.Range("I132").Value = 2218.2856
$ws.Range("K132").Value = 6654.8568
$ws.Range("M132").Value = -4124.8568
$ws.Range("H133").Value = 93833
$ws.Range("J133").Value = 93833
$ws.Range("L133").Value = 93833
$ws.Range("N133").Value = -103953
